$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D27").Value = "더 나은 생성모델을 위해 RLHF로 피드백 학습시키기"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/luda-rlhf/"

$ws.Range("D32").Value = "어트리뷰트 (Attribute), 프로퍼티 (Property), 디스크립터 (Descriptor)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/453"

$ws.Range("D51").Value = "[pandas] 데이터프레임 열 또는 행 기준으로 요소 총합 구하기, sum() 메서드"
$ws.Range("E51").Value = "https://bskyvision.com/entry/pandas-%EB%8D%B0%EC%9D%B4%ED%84%B0%ED%94%84%EB%A0%88%EC%9E%84-%EC%97%B4-%EB%98%90%EB%8A%94-%ED%96%89-%EA%B8%B0%EC%A4%80%EC%9C%BC%EB%A1%9C-%EC%9A%94%EC%86%8C-%EC%B4%9D%ED%95%A9-%EA%B5%AC%ED%95%98%EA%B8%B0-sum-%EB%A9%94%EC%84%9C%EB%93%9C"

$wb.Save()
